$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.834.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.398.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.52%  "

# Row 9
$ws.Range("E9").Value = "  -1.85%  "

# Row 10
$ws.Range("E10").Value = "  -1.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000170"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.829.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.758.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.395.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "320.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "562.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.89%  "

# Row 27
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.519.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0925"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.43%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.41%  "

# Row 31
$ws.Range("E31").Value = "  -5.48%  "

# Row 32
$ws.Range("E32").Value = "  -1.85%  "

# Row 33
$ws.Range("E33").Value = "  -0.57%  "

# Row 34
$ws.Range("E34").Value = "  -4.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.78%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.01%  "

# Row 38
$ws.Range("E38").Value = "  -5.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.14%  "

# Row 41
$ws.Range("E41").Value = "  -6.49%  "

# Row 42
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.35%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.68%  "

# Row 46
$ws.Range("E46").Value = "  -3.30%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0915"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
